# Append the April-2025 trail income figures beneath the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrailIncome")

$names = @(
    "HARSHAD NAGTILAK",
    "HRUTWIK GARDI",
    "PRATIK RAUL",
    "PRATIK SHIRBHATE",
    "SARANG THAKREY",
    "SHUBHAM MUNDADA",
    "YUKTA SONIGRA"
)

$incomes = @(
    22913.3785,
    592.861,
    29943.4216,
    1752.8626499999998,
    557.0264,
    89968.73989999996,
    16846.567150000003
)

$startRow = 9
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $incomes[$i]
    $ws.Cells.Item($row, 3).Value = "2025-04"
}
